# Applies the changes described by the publications.xlsx diff:
#   - dissertations sheet: add two new dissertation rows (Owen LLW 2021, Ziman K 2022)
#   - talks sheet: remove now-unused empty "links_html" placeholder cells
#   - posters sheet: remove now-unused empty "links_html" placeholder cells and
#                    add two new Wetterhahn poster rows (Carstensen et al. 2024, Jha et al. 2023)

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) dissertations sheet: append two dissertation rows (rows 4 and 5)
# -----------------------------------------------------------------
$diss = $wb.Worksheets.Item("dissertations")

$diss.Range("B4").Value = "Modeling the fast-timescale network dynamics that underlie complex thought"
$diss.Range("D4").Value = 'Owen LLW (2021) Modeling the fast-timescale network dynamics that underlie complex thought. <em>Doctoral dissertation: Dartmouth College</em>, Hanover, NH.'

$diss.Range("B5").Value = "Attending and remembering the external world"
$diss.Range("D5").Value = 'Ziman K (2022) Attending and remembering the external world. <em>Doctoral dissertation: Dartmouth College</em>, Hanover, NH.'

# -----------------------------------------------------------------
# 2) talks sheet: clear the empty placeholder cells in column E that
#    don't carry any links_html content (rows 2,3,4,5,7,9,12,13)
# -----------------------------------------------------------------
$talks = $wb.Worksheets.Item("talks")

$talksEmptyRows = @(2, 3, 4, 5, 7, 9, 12, 13)
foreach ($r in $talksEmptyRows) {
    $talks.Range("E$r").ClearContents()
}

# -----------------------------------------------------------------
# 3) posters sheet: clear all the empty placeholder cells in column E
#    (rows 2 through 43), then append two new poster rows (44 and 45)
# -----------------------------------------------------------------
$posters = $wb.Worksheets.Item("posters")

$posters.Range("E2:E43").ClearContents()

$posters.Range("B44").Value = "Translating neurophysiological recordings into dynamic estimates of conceptual knowledge and learning"
$posters.Range("C44").Value = "https://digitalcommons.dartmouth.edu/wetterhahn_2024/15/"
$posters.Range("D44").Value = 'Carstensen DL, Manning JR, Mucha P (2024) Translating neurophysiological recordings into dynamic estimates of conceptual knowledge and learning. <em>Wetterhahn Science Symposium</em>, Hanover, NH.'

$posters.Range("B45").Value = "Exploring high-order network dynamics in brains and stock markets"
$posters.Range("C45").Value = "https://digitalcommons.dartmouth.edu/wetterhahn_2023/4/"
$posters.Range("D45").Value = 'Jha K, Carstensen DL, Patel A, Manning JR (2023) Exploring high-order network dynamics in brains and stock markets. <em>Wetterhahn Science Symposium</em>, Hanover, NH.'
